$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest cryptos data refresh.
# Price values are forced to remain plain text (matching the sheet's existing
# inline-string convention) by briefly applying a text number format, then
# clearing the cell format afterwards so no residual style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.883.18'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.51'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.97%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3869'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3811'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.80'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.348'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08450'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.91'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.073'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.050'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("E16").Value = '  -2.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.650.99'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.07'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06990'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.61'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.50%  '

$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.79'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.908.34'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.974'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.08'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.408'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.90'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.845'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.501'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.832.26'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.65%  '

$ws.Range("E34").Value = '  -3.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08138'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.694'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02919'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.80'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2677'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09114'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7579'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.50'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.425'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.29'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6932'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.461'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.098'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.80%  '

$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("E49").Value = '  -1.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.33'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.234'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.64%  '

